$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Cell B4 holds the Cypher query that drives the "Case Files" tab export.
# Fix the column ordering in its RETURN clause: `Format` should immediately
# follow `File Name`, and the `Size` expression should immediately follow
# `File Type`, ahead of `Association`/`Description` (previously `Format`
# and `Size` were declared last, after `Description`).
$cell = $ws.Cells.Item(4, 2)
$lines = $cell.Text -split "`n"

# Lines (0-based) inside the RETURN clause, as originally ordered:
#  22 File Name
#  23 File Type
#  24 Association
#  25 Description
#  26 Format          <- move up, right after File Name
#  27 Size (CASE ...)  <- move up, right after File Type
#  28-31 Sample ID / Case ID / Breed / Diagnosis
$formatLine = "       " + $lines[26].TrimStart()
$sizeLine   = "       " + $lines[27].TrimStart()

$newLines = $lines[0..21] + $lines[22] + $formatLine + $lines[23] + $sizeLine + $lines[24] + $lines[25] + $lines[28..31]

$cell.Value = $newLines -join "`n"

# The saved cursor/selection on this sheet moved from B4 to C4.
$ws.Range("C4").Select()
